$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "43.124.37"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -1.46%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.279.98"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  -0.29%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "112.87"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.84%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "264.99"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.59%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.616"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.56%  "
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("E9").Value = "  -2.59%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "47.49"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.00%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0930"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.52%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "8.93"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.77%  "
$ws.Range("E13").Value = "  +0.95%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "15.54"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.01%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "2.621.78"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.13%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.861"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.20%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.285.38"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -1.03%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "43.162.56"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("E19").Value = "  -2.38%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "6.79"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +3.37%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "71.45"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.73%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "2.53"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.66%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "9.61"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.27%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.85"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.17%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.02%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "11.33"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -2.40%  "
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "40.45"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -8.77%  "
$ws.Range("B29").Value = "WEMIXToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "3.34"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -2.21%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "2.25"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -1.02%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "171.66"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -3.56%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "21.40"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -2.21%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.0904"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.68%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.78"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +3.97%  "
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.127"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "4.64"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -1.97%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.90"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.0355"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.93%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.104"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -5.39%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.64"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +9.38%  "
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "78.48"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +5.53%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "13.85"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +2.90%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.238"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -3.19%  "
$ws.Range("B44").Value = "THORChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "6.23"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +4.56%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.39"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -4.85%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "8.68"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.59%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "103.74"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +2.11%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.25"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.31%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0994"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.57%  "
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.438"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -4.19%  "
